$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each element is the new set of values (columns B..K) for data rows 2..24.
# This reflects each row of forecast-error values shifting one column to the
# left (the former first value is dropped) now that the ifoCAST full series
# evaluation adds one more leading comparison column upstream.
$rowValues = @(
    @(10.2000313112857, -7.713854356263806, 0.4730551773759682, 1.112750016219461, -0.688116008571576, -1.647050705078689, 1.1866675198023, -0.4824371439038641, 0.5377352938393787),
    @(-8.028779524684778, 0.1581300089549967, 0.7978248477984897, -1.003041176992548, -1.96197587349966, 0.8717423513813284, -0.7973623123248356, 0.2228101254184071, -1.490979331480904),
    @(-0.898244951227533, -0.2585501123840399, -2.059416137175077, -3.01835083368219, -0.1846326088012012, -1.853737272507365, -0.8335648347641225, -2.547354291663434, 0.281001081935538),
    @(0.6120962635868906, -1.188769761204147, -2.14770445771126, 0.6860137671697293, -0.9830908965364347, 0.03708154120680809, -1.676707915692503, 1.151647457906469, -1.083685141959651),
    @(-1.19942097753403, -2.158355674041143, 0.6753625508398458, -0.9937421128663182, 0.02643032487692459, -1.687359132022387, 1.140996241576585, -1.094336358289534, 0.2569986996281282),
    @(-2.266970206257284, 0.5667480186237051, -1.102356645082459, -0.08218420733921622, -1.795973664238528, 1.032381709360444, -1.202950890505675, 0.1483841674119874, -0.4730537624049144),
    @(0.774849739591444, -0.89425492411472, 0.1259175136285228, -1.587871943270789, 1.240483430328183, -0.994849169537936, 0.3564858883797264, -0.2649520414371754, 0.2356086557536556),
    @(-0.7777567537409195, 0.2424156840023232, -1.471373772896988, 1.356981600701984, -0.8783509991641355, 0.4729840587535268, -0.148453871063375, 0.352106826127456, -0.3158012462186854),
    @(0.2879090979994584, -1.425880358899853, 1.402475014699119, -0.8328575851670005, 0.5184774727506619, -0.1029604570662399, 0.3976002401245912, -0.2703078322215502, 0.1586931430164528),
    @(-1.443434480259818, 1.384920893339154, -0.8504117065269649, 0.5009233513906975, -0.1205145784262043, 0.3800461187646267, -0.2878619535815147, 0.1411390216564884, 0.1988105702346985),
    @(1.473028212290161, -0.7623043875759586, 0.5890306703417038, -0.0324072594751981, 0.4681534377156329, -0.1997546346305085, 0.2292463406074946, 0.2869178891857047, 0.4108842600239903),
    @(-0.5837297540881751, 0.7676053038294873, 0.1461673740125855, 0.6467280712034165, -0.02118000114272489, 0.4078209740952782, 0.4654925226734883, 0.5894588935117738, -0.303959229723018),
    @(1.679632531582137, 1.058194601765235, 1.558755298956066, 0.8908472266099251, 1.319848201847928, 1.377519750426138, 1.501486121264424, 0.608067998029632, 1.378198724973394),
    @(0.1142203657994787, 0.6147810629903097, -0.0531270093558317, 0.3758739658821714, 0.4335455144603815, 0.557511885298667, -0.3359062379361248, 0.4342244890076376, 0.1683237681281231),
    @(0.6187489605034189, -0.04915911184272259, 0.3798418633952805, 0.4375134119734906, 0.5614797828117761, -0.3319383404230157, 0.4381923865207467, 0.1722916656412322),
    @(0.08648097832751878, 0.5154819535655218, 0.573153502143732, 0.6971198729820175, -0.1962982502527744, 0.5738324766909881, 0.3079317558114735),
    @(0.2746757717098572, 0.3323473202880673, 0.4563136911263528, -0.4371044321084391, 0.3330262948353234, 0.06712557395580883),
    @(0.2870161050359709, 0.4109824758742565, -0.4824356473605354, 0.287695079583227, 0.02179435870371246),
    @(0.3441210539382026, -0.5492970692965893, 0.2208336576471732, -0.04506706323234141),
    @(-0.5788832716533059, 0.1912474552904566, -0.07465326558905801),
    @(0.1730967985608157, -0.0928039223186989),
    @(-0.1108357465673982),
    @()
)

$startRow = 2
$startCol = 2  # column B
$lastCol = 11  # column K

for ($i = 0; $i -lt $rowValues.Count; $i++) {
    $r = $startRow + $i
    $vals = $rowValues[$i]
    for ($c = $startCol; $c -le $lastCol; $c++) {
        $idx = $c - $startCol
        if ($idx -lt $vals.Count) {
            $ws.Cells.Item($r, $c).Value = $vals[$idx]
        } else {
            $ws.Cells.Item($r, $c).ClearContents()
        }
    }
}